# Crayfish Guide — "Oscilloscope converted to App Designer"
#
# The paragraph describing the "Single Trace" button is rewritten to
# describe the new "Stop" / "Stopping" button behaviour instead. The
# edit also relocates the document's (hidden) "_GoBack" bookmark from
# the very last (empty) paragraph to the point right after the new
# "Stopping"" text, which is exactly what Word does automatically when
# you save after editing at a new location.

$d = $word.ActiveDocument

$dq_open  = [char]0x201C   # “
$dq_close = [char]0x201D   # ”

# --- locate the span of text to replace -----------------------------
# Start: the closing curly quote immediately before "button allows you
# to plot a new single spike" (".. the "Re-trigger" button allows ..").
$rngStart = $d.Content
[void]$rngStart.Find.Execute($dq_close + " button allows you to plot a new single spike", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rngStart.Start

# End: the end of the trailing sentence that is being removed entirely.
$rngEnd = $d.Content
[void]$rngEnd.Find.Execute("full trace showing.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$end = $rngEnd.End

# --- build the replacement text --------------------------------------
$part1 = $dq_close + " button allows you to plot a new single spike.  "
$part2 = "When you stop the program running with the"
$part3 = " " + $dq_open
$part4 = "Stop"
$part5 = $dq_close + " "
$part6 = "button, it will continue collecting data until the trace display is filled.  This allows you to easily save a complete trace. If you want the program to stop immediately, press the button while it displays "
$part7 = $dq_open
$part8 = "Stopping" + $dq_close
$part9 = "."

$newText = $part1 + $part2 + $part3 + $part4 + $part5 + $part6 + $part7 + $part8 + $part9

$full = $d.Range($start, $end)
$full.Text = $newText

# --- re-apply bold formatting to "Stop" and ""Stopping"" -------------
$stopStart = $start + $part1.Length + $part2.Length + $part3.Length
$stopEnd   = $stopStart + $part4.Length
$d.Range($stopStart, $stopEnd).Font.Bold = 1

$quoteStart   = $stopEnd + $part5.Length + $part6.Length
$quoteMid     = $quoteStart + $part7.Length
$quoteEnd     = $quoteMid + $part8.Length
$d.Range($quoteStart, $quoteMid).Font.Bold = 1
$d.Range($quoteMid, $quoteEnd).Font.Bold = 1

# --- move the "_GoBack" bookmark to just after ""Stopping"" ----------
# Adding a bookmark with a name that already exists elsewhere in the
# document relocates it (bookmark names are unique), which removes it
# from the empty trailing paragraph automatically.
$d.Bookmarks.Add("_GoBack", $d.Range($quoteEnd, $quoteEnd))
